# Apply cryptos list price/volume update (GitHub Actions refresh)
# Leading apostrophe forces text interpretation so numeric-looking strings
# (e.g. "4.50", "0.0910") keep their exact formatting instead of being
# auto-converted to numbers (which would strip trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.139.22"
$ws.Range("E2").Value = "'  -0.39%  "
$ws.Range("D3").Value = "'2.076.47"
$ws.Range("E3").Value = "'  -0.96%  "
$ws.Range("E4").Value = "'  +0.11%  "
$ws.Range("D5").Value = "'253.33"
$ws.Range("E5").Value = "'  +1.08%  "
$ws.Range("E6").Value = "'  +2.28%  "
$ws.Range("D7").Value = "'59.15"
$ws.Range("E7").Value = "'  +9.54%  "
$ws.Range("E8").Value = "'  +0.02%  "
$ws.Range("E9").Value = "'  +4.49%  "
$ws.Range("E10").Value = "'  -0.46%  "
$ws.Range("E11").Value = "'  +7.85%  "
$ws.Range("E12").Value = "'  +2.54%  "
$ws.Range("E13").Value = "'  +6.70%  "
$ws.Range("D14").Value = "'2.379.99"
$ws.Range("E14").Value = "'  -0.85%  "
$ws.Range("D15").Value = "'0.824"
$ws.Range("E15").Value = "'  -2.03%  "
$ws.Range("E16").Value = "'  +6.97%  "
$ws.Range("D17").Value = "'2.076.72"
$ws.Range("E17").Value = "'  -0.90%  "
$ws.Range("D18").Value = "'37.258.52"
$ws.Range("E18").Value = "'  -0.05%  "
$ws.Range("D19").Value = "'15.81"
$ws.Range("E19").Value = "'  +7.39%  "
$ws.Range("D20").Value = "'74.66"
$ws.Range("E20").Value = "'  +2.38%  "
$ws.Range("D21").Value = "'0.0₃0927"
$ws.Range("E21").Value = "'  +9.84%  "
$ws.Range("E22").Value = "'  +5.14%  "
$ws.Range("D23").Value = "'239.62"
$ws.Range("E23").Value = "'  -0.67%  "
$ws.Range("E24").Value = "'  -0.04%  "
$ws.Range("D25").Value = "'2.41"
$ws.Range("E25").Value = "'  -2.62%  "
$ws.Range("D26").Value = "'2.27"
$ws.Range("E26").Value = "'  +14.10%  "
$ws.Range("D27").Value = "'169.81"
$ws.Range("E27").Value = "'  -1.48%  "
$ws.Range("D28").Value = "'9.34"
$ws.Range("E28").Value = "'  +0.75%  "
$ws.Range("D29").Value = "'20.38"
$ws.Range("E29").Value = "'  -1.40%  "
$ws.Range("D30").Value = "'0.126"
$ws.Range("E30").Value = "'  +2.66%  "
$ws.Range("E31").Value = "'  +6.02%  "
$ws.Range("E32").Value = "'  +6.09%  "
$ws.Range("E33").Value = "'  +2.75%  "
$ws.Range("D34").Value = "'4.50"
$ws.Range("E34").Value = "'  +9.19%  "
$ws.Range("D35").Value = "'0.0910"
$ws.Range("E35").Value = "'  +0.45%  "
$ws.Range("E36").Value = "'  -0.07%  "
$ws.Range("E37").Value = "'  +0.83%  "
$ws.Range("B38").Value = "'WEMIXToken"
$ws.Range("C38").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "'1.77"
$ws.Range("E38").Value = "'  -3.63%  "
$ws.Range("B39").Value = "'Cronos"
$ws.Range("C39").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D39").Value = "'0.116"
$ws.Range("E39").Value = "'  +25.01%  "
$ws.Range("E40").Value = "'  +2.17%  "
$ws.Range("B41").Value = "'InjectiveProtocol"
$ws.Range("C41").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'17.92"
$ws.Range("E41").Value = "'  -2.81%  "
$ws.Range("B42").Value = "'VeChain"
$ws.Range("C42").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0227"
$ws.Range("E42").Value = "'  +0.59%  "
$ws.Range("E43").Value = "'  +0.39%  "
$ws.Range("D44").Value = "'99.19"
$ws.Range("E44").Value = "'  +0.22%  "
$ws.Range("D45").Value = "'4.33"
$ws.Range("E45").Value = "'  +7.41%  "
$ws.Range("E46").Value = "'  +1.53%  "
$ws.Range("D47").Value = "'4.57"
$ws.Range("E47").Value = "'  +13.82%  "
$ws.Range("D48").Value = "'2.50"
$ws.Range("E48").Value = "'  +8.47%  "
$ws.Range("D49").Value = "'1.307.16"
$ws.Range("E49").Value = "'  -0.93%  "
$ws.Range("E50").Value = "'  -0.26%  "
$ws.Range("E51").Value = "'  -0.86%  "
